# B1--and-B2-PowerPoint.pptx edit
# 1) Table on slide 5 switches to the built-in "No Style, Table Grid"
#    table style (GUID 03D99B5B-DB7C-46B8-84C0-AA6FFBEBFE95).
# 2) The deck's theme (ppt/theme/theme1.xml) switches from the
#    "Integral" / "Red Violet" color scheme to the stock PowerPoint
#    "Office Theme" / "Office" color scheme (font + format schemes are
#    already identical between the two themes, so only the 12 theme
#    colors need to change).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
$slide5 = $p.Slides.Item(5)
for ($i = 1; $i -le $slide5.Shapes.Count; $i++) {
    $shp = $slide5.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{03D99B5B-DB7C-46B8-84C0-AA6FFBEBFE95}")
    }
}

# --- 2) Theme colors --------------------------------------------------
# ThemeColorScheme index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
